$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# --- Update SecondSearchTerm (sheet2) selection ---
$sheet2.Range("B16").Select() | Out-Null

# --- New sheet: NumberOfFavorites ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$numFav = $wb.Worksheets.Add($null, $lastSheet)
$numFav.Name = "NumberOfFavorites"
$numFav.Range("A2").Value = "Queens, NY"
$sheet2.Range("A2").Copy()
$numFav.Range("B2").PasteSpecial($xlPasteFormats) | Out-Null
$numFav.Range("B2").Value = "2"
$numFav.Range("B5").Select() | Out-Null

# --- New sheet: OneFavorite ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oneFav = $wb.Worksheets.Add($null, $lastSheet)
$oneFav.Name = "OneFavorite"
$sheet2.Range("A2").Copy()
$oneFav.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null
$oneFav.Range("A2").Value = "Queens, NY"
$sheet2.Range("A2").Copy()
$oneFav.Range("B2").PasteSpecial($xlPasteFormats) | Out-Null
$oneFav.Range("B2").Value = "1"
$oneFav.Range("B5").Select() | Out-Null

# --- New sheet: RemoveFavorites ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$remFav = $wb.Worksheets.Add($null, $lastSheet)
$remFav.Name = "RemoveFavorites"
$remFav.Range("A2").Value = "Queens, NY"
$sheet1.Range("B2").Copy()
$remFav.Range("B2").PasteSpecial($xlPasteFormats) | Out-Null
$remFav.Range("B2").Value = "No Favorites Yet"
$remFav.Range("C7").Select() | Out-Null

$excel.CutCopyMode = $false
